$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.223.79'
$ws.Range('E2').Value = '  +2.42%  '
$ws.Range('D3').Value = '2.423.47'
$ws.Range('E3').Value = '  +3.53%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.23'
$ws.Range('E5').Value = '  +2.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.02'
$ws.Range('E6').Value = '  +4.94%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.62%  '
$ws.Range('D9').Value = '2.422.29'
$ws.Range('E9').Value = '  +3.48%  '
$ws.Range('E10').Value = '  +5.20%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.39'
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('E13').Value = '  +2.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.27'
$ws.Range('E14').Value = '  +7.06%  '
$ws.Range('E15').Value = '  +9.56%  '
$ws.Range('D16').Value = '2.862.26'
$ws.Range('E16').Value = '  +3.31%  '
$ws.Range('D17').Value = '62.115.62'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('D18').Value = '2.425.68'
$ws.Range('E18').Value = '  +3.31%  '
$ws.Range('E19').Value = '  +4.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '323.73'
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('E22').Value = '  +3.43%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  +5.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.93'
$ws.Range('E25').Value = '  +2.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.12'
$ws.Range('E26').Value = '  +7.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '568.63'
$ws.Range('E27').Value = '  +14.62%  '
$ws.Range('D28').Value = '2.545.40'
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.40'
$ws.Range('E30').Value = '  +5.75%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0939'
$ws.Range('E31').Value = '  +9.36%  '
$ws.Range('E32').Value = '  +6.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.148'
$ws.Range('E33').Value = '  +1.80%  '
$ws.Range('E34').Value = '  +4.16%  '
$ws.Range('E35').Value = '  +5.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.71'
$ws.Range('E36').Value = '  +9.00%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.82'
$ws.Range('E38').Value = '  +4.92%  '
$ws.Range('E39').Value = '  +2.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.88'
$ws.Range('E40').Value = '  +3.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.78'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '149.62'
$ws.Range('E42').Value = '  +5.34%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  +2.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.32'
$ws.Range('E45').Value = '  +14.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.77'
$ws.Range('E46').Value = '  +6.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.63'
$ws.Range('E47').Value = '  +2.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0543'
$ws.Range('E48').Value = '  +5.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.35'
$ws.Range('E49').Value = '  +7.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.589'
$ws.Range('E50').Value = '  +3.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0915'
$ws.Range('E51').Value = '  +1.58%  '
